$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("meta")
$wsEvents = $wb.Worksheets.Item("events")

# "meta" sheet: update the impacto_clave cell (C2) with the expanded text,
# and grow row 2 so the longer text is readable.
$wsMeta.Range("C2").Value = "PN Cuna Más: 01 CIAI afectado (28 niños SCD) - Ucayali - Padre Abad"
$wsMeta.Rows.Item(2).RowHeight = 30
$wsMeta.Range("C3").Select()

# "events" sheet: update the afectacion_midis cell (H5) for the
# Ucayali / Padre Abad row, then make "events" the active tab/selection.
$wsEvents.Activate()
$wsEvents.Range("H5").Value = "PN Cuna Más: 1 CIAI afectado (28 niños SCD)"
$wsEvents.Range("H6").Select()
